$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -7
$ws.Range("F3").Value = -3
$ws.Range("F4").Value = -9
$ws.Range("F8").Value = -5
